$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$src = $ws.Range("B64")
$src.Copy() | Out-Null
$dst = $ws.Range("G2")
$dst.PasteSpecial(-4122) | Out-Null   # xlPasteFormats = -4122
$dst.Value = "Thematic Area"
$dst.Font.Name = "Arial"
$dst.Font.Bold = $true
